$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8854930400848389
$ws.Range("B1").Value = 0.8070028424263
$ws.Range("C1").Value = 2.394325017929077
$ws.Range("D1").Value = 5.37716817855835
$ws.Range("E1").Value = 1.201570272445679
